$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 22-26 (the removed Resolving-Mac x Resolving-Mac block of rows)
$ws.Range("A22:T26").Delete() | Out-Null

# Update the numeric columns (G,H,I,J,M,N,O,P,Q,R,S,T) for rows 2-21 with the
# recomputed TPM-based values. Text columns (A-D) and count columns (E,F,K,L)
# are unchanged by this edit.

# Row 2
$ws.Range("G2").Value = 37.393558
$ws.Range("H2").Value = 112.180674
$ws.Range("I2").Value = 0.6214236533709717
$ws.Range("J2").Value = 0.6390053627425325
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 4559.349373663362
$ws.Range("R2").Value = 41034.14436297026
$ws.Range("S2").Value = 0.1418233041431385
$ws.Range("T2").Value = 0.1546214331950799

# Row 3
$ws.Range("G3").Value = 37.393558
$ws.Range("H3").Value = 112.180674
$ws.Range("I3").Value = 0.6214236533709717
$ws.Range("J3").Value = 0.6390053627425325
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 5531.043077886139
$ws.Range("R3").Value = 49779.38770097525
$ws.Range("S3").Value = 0.1720488474068329
$ws.Range("T3").Value = 0.1875745282225055

# Row 4
$ws.Range("G4").Value = 37.393558
$ws.Range("H4").Value = 112.180674
$ws.Range("I4").Value = 0.6214236533709717
$ws.Range("J4").Value = 0.6390053627425325
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 3122.547914054221
$ws.Range("R4").Value = 28102.93122648799
$ws.Range("S4").Value = 0.09713010042058083
$ws.Range("T4").Value = 0.1058951166322749

# Row 5
$ws.Range("G5").Value = 37.393558
$ws.Range("H5").Value = 112.180674
$ws.Range("I5").Value = 0.6214236533709717
$ws.Range("J5").Value = 0.6390053627425325
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 3405.373834464575
$ws.Range("R5").Value = 20432.24300678745
$ws.Range("S5").Value = 0.1059276948233304
$ws.Range("T5").Value = 0.07699107039138327

# Row 6
$ws.Range("G6").Value = 37.393558
$ws.Range("H6").Value = 112.180674
$ws.Range("I6").Value = 0.6214236533709717
$ws.Range("J6").Value = 0.6390053627425325
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 3359.27384087154
$ws.Range("R6").Value = 30233.46456784386
$ws.Range("S6").Value = 0.1044937065770889
$ws.Range("T6").Value = 0.1139232143012889

# Row 7
$ws.Range("H7").Value = 52.61241699999999
$ws.Range("I7").Value = 0.2914459257466844
$ws.Range("J7").Value = 0.299691697429509
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 2138.321886850721
$ws.Range("R7").Value = 19244.89698165649
$ws.Range("S7").Value = 0.06651472621653735
$ws.Range("T7").Value = 0.07251701233669881

# Row 8
$ws.Range("H8").Value = 52.61241699999999
$ws.Range("I8").Value = 0.2914459257466844
$ws.Range("J8").Value = 0.299691697429509
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 2594.043470078536
$ws.Range("R8").Value = 23346.39123070682
$ws.Range("S8").Value = 0.08069042002847709
$ws.Range("T8").Value = 0.08797192016711122

# Row 9
$ws.Range("H9").Value = 52.61241699999999
$ws.Range("I9").Value = 0.2914459257466844
$ws.Range("J9").Value = 0.299691697429509
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 1464.466089379182
$ws.Range("R9").Value = 13180.19480441263
$ws.Range("S9").Value = 0.04555374080369202
$ws.Range("T9").Value = 0.04966450847425716

# Row 10
$ws.Range("H10").Value = 52.61241699999999
$ws.Range("I10").Value = 0.2914459257466844
$ws.Range("J10").Value = 0.299691697429509
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 1597.110641532954
$ws.Range("R10").Value = 9582.663849197725
$ws.Range("S10").Value = 0.04967978755319121
$ws.Range("T10").Value = 0.03610859300691854

# Row 11
$ws.Range("H11").Value = 52.61241699999999
$ws.Range("I11").Value = 0.2914459257466844
$ws.Range("J11").Value = 0.299691697429509
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 1575.489875672569
$ws.Range("R11").Value = 14179.40888105313
$ws.Range("S11").Value = 0.04900725114478671
$ws.Range("T11").Value = 0.05342966344452321

# Row 12
$ws.Range("G12").Value = 0.276071
$ws.Range("H12").Value = 0.8282130000000001
$ws.Range("I12").Value = 0.004587877126048758
$ws.Range("J12").Value = 0.004717680235127497
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 33.660988904469
$ws.Range("R12").Value = 302.948900140221
$ws.Range("S12").Value = 0.001047060068424096
$ws.Range("T12").Value = 0.001141546725336993

# Row 13
$ws.Range("G13").Value = 0.276071
$ws.Range("H13").Value = 0.8282130000000001
$ws.Range("I13").Value = 0.004587877126048758
$ws.Range("J13").Value = 0.004717680235127497
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 40.83485699742999
$ws.Range("R13").Value = 367.51371297687
$ws.Range("S13").Value = 0.001270210696517613
$ws.Range("T13").Value = 0.00138483445680444

# Row 14
$ws.Range("G14").Value = 0.276071
$ws.Range("H14").Value = 0.8282130000000001
$ws.Range("I14").Value = 0.004587877126048758
$ws.Range("J14").Value = 0.004717680235127497
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 23.05330038882267
$ws.Range("R14").Value = 207.479703499404
$ws.Range("S14").Value = 0.0007170968847952413
$ws.Range("T14").Value = 0.0007818076777767111

# Row 15
$ws.Range("G15").Value = 0.276071
$ws.Range("H15").Value = 0.8282130000000001
$ws.Range("I15").Value = 0.004587877126048758
$ws.Range("J15").Value = 0.004717680235127497
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 25.1413615108375
$ws.Range("R15").Value = 150.848169065025
$ws.Range("S15").Value = 0.0007820481976486875
$ws.Range("T15").Value = 0.0005684134629290845

# Row 16
$ws.Range("G16").Value = 0.276071
$ws.Range("H16").Value = 0.8282130000000001
$ws.Range("I16").Value = 0.004587877126048758
$ws.Range("J16").Value = 0.004717680235127497
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 24.80101220973
$ws.Range("R16").Value = 223.20910988757
$ws.Range("S16").Value = 0.0007714612786631195
$ws.Range("T16").Value = 0.0008410779122802684

# Row 17
$ws.Range("G17").Value = 4.9669165
$ws.Range("H17").Value = 9.933833
$ws.Range("I17").Value = 0.08254254375629515
$ws.Range("J17").Value = 0.05658525959283094
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 605.6098655632935
$ws.Range("R17").Value = 3633.659193379761
$ws.Range("S17").Value = 0.01883812472279512
$ws.Range("T17").Value = 0.01369205087482877

# Row 18
$ws.Range("G18").Value = 4.9669165
$ws.Range("H18").Value = 9.933833
$ws.Range("I18").Value = 0.08254254375629515
$ws.Range("J18").Value = 0.05658525959283094
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 734.6781262634449
$ws.Range("R18").Value = 4408.06875758067
$ws.Range("S18").Value = 0.02285292720716709
$ws.Range("T18").Value = 0.01661011627026021

# Row 19
$ws.Range("G19").Value = 4.9669165
$ws.Range("H19").Value = 9.933833
$ws.Range("I19").Value = 0.08254254375629515
$ws.Range("J19").Value = 0.05658525959283094
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 414.7622100137273
$ws.Range("R19").Value = 2488.573260082364
$ws.Range("S19").Value = 0.01290160990900197
$ws.Range("T19").Value = 0.009377233766134628

# Row 20
$ws.Range("G20").Value = 4.9669165
$ws.Range("H20").Value = 9.933833
$ws.Range("I20").Value = 0.08254254375629515
$ws.Range("J20").Value = 0.05658525959283094
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 452.3294490208813
$ws.Range("R20").Value = 1809.317796083525
$ws.Range("S20").Value = 0.01407017794950041
$ws.Range("T20").Value = 0.00681772009819843

# Row 21
$ws.Range("G21").Value = 4.9669165
$ws.Range("H21").Value = 9.933833
$ws.Range("I21").Value = 0.08254254375629515
$ws.Range("J21").Value = 0.05658525959283094
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 446.206072934895
$ws.Range("R21").Value = 2677.23643760937
$ws.Range("S21").Value = 0.01387970396783055
$ws.Range("T21").Value = 0.0100881385834089
